$d = $word.ActiveDocument

$replacements = @(
    @("73×24=1752", "58×38=2204"),
    @("67×21=1407", "34×67=2278"),
    @("84×21=1764", "77×84=6468"),
    @("20×70=1400", "99×39=3861"),
    @("91×33=3003", "12×94=1128"),
    @("94×13=1222", "80×59=4720"),
    @("83×14=1162", "78×13=1014"),
    @("44×53=2332", "34×87=2958"),
    @("82×82=6724", "47×98=4606"),
    @("70×83=5810", "89×97=8633"),
    @("92×68=6256", "38×63=2394"),
    @("61×15=915",  "63×14=882"),
    @("11×64=704",  "57×81=4617"),
    @("37×65=2405", "65×91=5915"),
    @("52×53=2756", "12×18=216"),
    @("40×51=2040", "86×43=3698"),
    @("98×90=8820", "91×60=5460"),
    @("85×15=1275", "30×17=510"),
    @("87×11=957",  "96×31=2976"),
    @("91×42=3822", "16×51=816"),
    @("96×40=3840", "22×68=1496"),
    @("65×60=3900", "80×55=4400"),
    @("15×14=210",  "99×29=2871"),
    @("54×80=4320", "33×83=2739"),
    @("95×68=6460", "63×16=1008")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}

$d.Save()
